$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need text (General->Text) formatting to avoid Excel
# auto-converting numeric-looking / percent-looking strings into numbers.
$textCells = @(
    "D2",
    "E2",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)

foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New cell values as they appear in the updated symbol list.
$updates = @{
    "D2" = "313.85"
    "E2" = "-0.59%"
    "E3" = "-3.21%"
    "D4" = "5.090"
    "E4" = "-0.72%"
    "D5" = "0.07751"
    "E5" = "-5.47%"
    "B6" = "GateToken"
    "C6" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D6" = "4.357"
    "E6" = "-0.19%"
    "B7" = "FTXToken"
    "C7" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D7" = "1.918"
    "E7" = "-3.90%"
    "D8" = "8.208"
    "E8" = "-1.34%"
    "D9" = "0.9181"
    "E9" = "-2.20%"
    "D10" = "0.1254"
    "E10" = "-3.81%"
    "D11" = "0.1894"
    "E11" = "-3.90%"
    "D12" = "0.08881"
    "E12" = "-2.34%"
    "D13" = "0.03430"
    "E13" = "-1.46%"
    "E14" = "-0.51%"
    "D15" = "0.001369"
    "E15" = "-2.97%"
    "D16" = "0.006048"
    "E16" = "-7.82%"
    "D17" = "3.531"
    "E17" = "-2.81%"
    "D18" = "3.098"
    "E18" = "-6.12%"
    "E19" = "-1.83%"
    "B20" = "ProBitToken"
    "C20" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "D20" = "0.1280"
    "E20" = "-2.68%"
    "B21" = "MCDex"
    "C21" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D21" = "5.027"
    "E21" = "1.47%"
    "D22" = "0.2592"
    "E22" = "4.19%"
    "D23" = "0.02105"
    "E23" = "5,591.90%"
    "D24" = "0.04400"
    "E24" = "1.08%"
    "D25" = "0.001210"
    "E25" = "-2.32%"
    "D26" = "0.004246"
    "E26" = "-11.06%"
    "D27" = "0.0001351"
    "E27" = "-65.26%"
    "D39" = "0.02137"
    "E39" = "-4.64%"
    "D40" = "0.04991"
    "E40" = "-4.43%"
    "D41" = "0.007865"
    "E41" = "1.43%"
    "D42" = "0.009974"
    "E42" = "-3.41%"
    "D43" = "0.1345"
    "E43" = "-3.98%"
    "D44" = "0.002061"
    "E44" = "-1.81%"
    "D45" = "0.009681"
    "E45" = "-0.69%"
    "D46" = "0.00006497"
    "E46" = "-4.28%"
    "D47" = "0.00000000750"
    "E47" = "0.09%"
    "D48" = "0.003199"
    "E48" = "11.03%"
    "E49" = "-0.04%"
    "D50" = "0.00002101"
    "E50" = "0.09%"
    "D51" = "0.0002001"
    "E51" = "0.09%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
